$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Range("H9").Formula = "=455/0.89"
$ws.Range("I9").Value = 445

# Row 10
$ws.Range("H10").Value = 1
$ws.Range("I10").Formula = "=C9"
$ws.Range("J10").Formula = "=I10+H10"

# Row 11
$ws.Range("H11").Formula = "=H9*H10"
$ws.Range("I11").Formula = "=I9*I10"
$ws.Range("J11").Formula = "=SUM(H11:I11)/J10"

# Row 13
$ws.Range("H13").Formula = "=455/0.86"
$ws.Range("I13").Value = 445

# Row 14
$ws.Range("H14").Value = 1
$ws.Range("I14").Formula = "=C9"
$ws.Range("J14").Formula = "=I14+H14"

# Row 15
$ws.Range("H15").Formula = "=H13*H14"
$ws.Range("I15").Formula = "=I13*I14"
$ws.Range("J15").Formula = "=SUM(H15:I15)/J14"

# Update selection to match target state
$ws.Range("I23").Select()
